$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Constant used for PasteSpecial(xlPasteFormats)
$xlPasteFormats = -4122

# ------------------------------------------------------------------
# 1) Creditos-aula: "4" -> "2"  (row 5)
#    Use a leading apostrophe so Excel stores it as TEXT (matching the
#    original file, where these look-like-number values are text), then
#    restore the original (non quote-prefixed) cell formatting by
#    pasting formats from a sibling cell that already carries the
#    correct column style.
# ------------------------------------------------------------------
$ws.Range("B5").Value = "'2"
$ws.Range("C5").Value = "'2"
$ws.Range("B6").Copy()
$ws.Range("B5").PasteSpecial($xlPasteFormats)
$ws.Range("C6").Copy()
$ws.Range("C5").PasteSpecial($xlPasteFormats)

# ------------------------------------------------------------------
# 2) Carga horaria: "60 h" -> "30 h"  (row 7)
# ------------------------------------------------------------------
$ws.Range("B7").Value = "'30 h"
$ws.Range("C7").Value = "'30 h"
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial($xlPasteFormats)
$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial($xlPasteFormats)

# ------------------------------------------------------------------
# 3) Ativacao: "01/01/2012" -> "01/01/2023" (row 8, and also mirrored
#    at row 15 which shares the same underlying text in the source
#    workbook).
# ------------------------------------------------------------------
$ws.Range("B8").Value = "'01/01/2023"
$ws.Range("C8").Value = "'01/01/2023"
$ws.Range("B9").Copy()
$ws.Range("B8").PasteSpecial($xlPasteFormats)
$ws.Range("C9").Copy()
$ws.Range("C8").PasteSpecial($xlPasteFormats)

$ws.Range("B15").Value = "'01/01/2023"
$ws.Range("C15").Value = "'01/01/2023"
$ws.Range("B9").Copy()
$ws.Range("B15").PasteSpecial($xlPasteFormats)
$ws.Range("C9").Copy()
$ws.Range("C15").PasteSpecial($xlPasteFormats)

# ------------------------------------------------------------------
# 4) New English "Objectives" paragraph (row 11). B11/C11 were empty
#    before, so after writing the value we copy formats from a
#    well-formed sibling row (13) to make sure the new cells pick up
#    the correct column styles (s=2 for B, s=3 for C) instead of the
#    engine's default column-style resolution.
# ------------------------------------------------------------------
$objectivesText = "Provide students with basic knowledge of optical metrology, ie methods of measuring the size and geometry of mechanical components using optical methods, with particular emphasis on laser interferometry."
$ws.Range("B11").Value = $objectivesText
$ws.Range("C11").Value = $objectivesText
$ws.Range("B13").Copy()
$ws.Range("B11").PasteSpecial($xlPasteFormats)
$ws.Range("C13").Copy()
$ws.Range("C11").PasteSpecial($xlPasteFormats)

# ------------------------------------------------------------------
# 5) New English "Short syllabus" paragraph (row 14).
# ------------------------------------------------------------------
$shortSyllabusText = "Present the main optical techniques for measuring quantities such as length, displacement and shape, with emphasis on laser interferometric techniques."
$ws.Range("B14").Value = $shortSyllabusText
$ws.Range("C14").Value = $shortSyllabusText
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial($xlPasteFormats)
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial($xlPasteFormats)

# ------------------------------------------------------------------
# 6) New English "Syllabus" paragraph (row 16).
# ------------------------------------------------------------------
$syllabusText = 'Electromagnetic theory of light: notions of mathematical representation of the light wave and interpretation of phenomena such as polarization, interference and diffraction. Refraction, reflection and geometric optics: Snell''s laws, Fraunhofer equations, total reflection and geometric optics. Propagation of light in special media such as fiber optic crystals. Fourier optics and holography: Fourier transform and its application in optics as a case of special filters and halography. Light sources and sensors: definition and description of incoherent and coherent sources and description of point, position and image sensors. Optical components and tuning of optical systems. Length measurement: method such as interferometry, Moirè fringes, methods for measuring large distances. Shape measurement: various methods and techniques for geometric shape measurement. Displacement, deformation and vibration measurement: measurement methods employing holography, speckle" and Moirè fringes. Velocity measurement: speed measurement methods and optical fiber sensors. Fault inspection: methods for inspecting geometric and internal spaces using diffraction or scattering of light.'
$ws.Range("B16").Value = $syllabusText
$ws.Range("C16").Value = $syllabusText
$ws.Range("B13").Copy()
$ws.Range("B16").PasteSpecial($xlPasteFormats)
$ws.Range("C13").Copy()
$ws.Range("C16").PasteSpecial($xlPasteFormats)
